$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Strategy column (L) values from "MACD" to "EarlyMACD" for the two data rows
$ws.Range("L2").Value = "EarlyMACD"
$ws.Range("L3").Value = "EarlyMACD"

# Move the active selection to L3
$ws.Activate()
$ws.Range("L3").Select()
